$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a blank row above row 1 and a blank column before column A on
#    both worksheets (the lab tables grew a spacer row/column so the header
#    block could be widened).
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Table 1a")
$ws2 = $wb.Worksheets.Item("Table 1b")

foreach ($ws in @($ws1, $ws2)) {
    $ws.Rows.Item(1).Insert()
    $ws.Columns.Item(1).Insert()
    $ws.Rows.Item(1).RowHeight = 17
}

# ---------------------------------------------------------------------------
# 2. Fill in the measured/derived data that used to be blank placeholder
#    cells. Column layout after the insert (both sheets):
#      B = case label (1a/1b/2a/2b)      -- already populated
#      C = fast-side flow rate (kg/s)
#      D = "(fast)"                      -- already populated
#      E = slow-side flow rate (kg/s)
#      F = "(slow)"                      -- already populated
#    Table 1a then has Temperature/Ui/Heat-transfer columns G:L; Table 1b has
#    the NTU/Cr/effectiveness columns G:K.
# ---------------------------------------------------------------------------

function Set-Col($ws, $rng, $values) {
    $arr = New-Object 'object[,]' $values.Length, 1
    for ($i = 0; $i -lt $values.Length; $i++) { $arr[$i, 0] = $values[$i] }
    $ws.Range($rng).Value = $arr
}

# -- Table 1a ---------------------------------------------------------------
Set-Col $ws1 "C4:C7" @(0.1766, 0.172, 0.1487, 0.1581)
Set-Col $ws1 "E4:E7" @(0.1962, 0.128, 0.2214, 0.1278)
Set-Col $ws1 "G4:G7" @(5.166, 7, 4, 6.89)
Set-Col $ws1 "H4:H7" @(6.78, 6.78, 7.17, 7)
Set-Col $ws1 "I4:I7" @(5290, 4670, 4580, 4530)
Set-Col $ws1 "J4:J7" @(-6.54, -6.37, -5.82, -6.04)
Set-Col $ws1 "K4:K7" @(6.58, 5.72, 5.67, 5.64)
Set-Col $ws1 "L4:L7" @(0.95, 10.64, 2.57, 6.9)

# -- Table 1b ----------------------------------------------------------------
Set-Col $ws2 "C4:C7" @(0.1766, 0.172, 0.1487, 0.1581)
Set-Col $ws2 "E4:E7" @(0.1962, 0.128, 0.2214, 0.1278)
Set-Col $ws2 "G4:G7" @(0.7695, 0.8705, 0.5727, 0.9483)
Set-Col $ws2 "H4:H7" @(0.2759, 0.243, 0.2387, 0.2357)
Set-Col $ws2 "I4:I7" @(0.2283, 0.1961, 0.2001, 0.1905)
Set-Col $ws2 "J4:J7" @(0.2217, 0.198, 0.2008, 0.1917)
Set-Col $ws2 "K4:K7" @(2.9816, 0.9723, 0.3848, 0.6449)

# ---------------------------------------------------------------------------
# 3. View/window state: Table 1b becomes the active sheet/tab, each sheet
#    gets its own zoom level and the analysis selections are left on the
#    ranges the author was last looking at.
# ---------------------------------------------------------------------------
$ws1.Select()
$excel.ActiveWindow.Zoom = 170
$ws1.Range("C4:F7").Select()

$ws2.Select()
$excel.ActiveWindow.Zoom = 190
$ws2.Range("G15").Select()
